$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 290, shifting the existing
# rows 290-294 down to 293-297.
$ws.Range("A290:A292").EntireRow.Insert()

# Row 290 - new "1a nueva(o)" record
$ws.Range("A290").Value = 9
$ws.Range("B290").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C290").Value = "Metropolitana"
$ws.Range("D290").Value = 44890
$ws.Range("E290").Value = 13
$ws.Range("F290").Value = 100112003
$ws.Range("G290").Value = "Ajo"
$ws.Range("H290").Value = "Chino"
$ws.Range("I290").Value = "1a nueva(o)"
$ws.Range("J290").Value = 950
$ws.Range("K290").Value = 1500
$ws.Range("L290").Value = 1500
$ws.Range("M290").Value = 1500
$ws.Range("N290").Value = "`$/paquete 20 unidades (volumen en unidades)"
$ws.Range("O290").Value = "Provincia de Talagante"
$ws.Range("P290").Value = 75
$ws.Range("Q290").Value = 20
$ws.Range("R290").Value = "Hortaliza"

# Row 291 - new "2a nueva(o)" record
$ws.Range("A291").Value = 9
$ws.Range("B291").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C291").Value = "Metropolitana"
$ws.Range("D291").Value = 44890
$ws.Range("E291").Value = 13
$ws.Range("F291").Value = 100112003
$ws.Range("G291").Value = "Ajo"
$ws.Range("H291").Value = "Chino"
$ws.Range("I291").Value = "2a nueva(o)"
$ws.Range("J291").Value = 600
$ws.Range("K291").Value = 1200
$ws.Range("L291").Value = 1200
$ws.Range("M291").Value = 1200
$ws.Range("N291").Value = "`$/paquete 20 unidades (volumen en unidades)"
$ws.Range("O291").Value = "Provincia de Talagante"
$ws.Range("P291").Value = 60
$ws.Range("Q291").Value = 20
$ws.Range("R291").Value = "Hortaliza"

# Row 292 - new "3a nueva (o)" record
$ws.Range("A292").Value = 9
$ws.Range("B292").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C292").Value = "Metropolitana"
$ws.Range("D292").Value = 44890
$ws.Range("E292").Value = 13
$ws.Range("F292").Value = 100112003
$ws.Range("G292").Value = "Ajo"
$ws.Range("H292").Value = "Chino"
$ws.Range("I292").Value = "3a nueva (o)"
$ws.Range("J292").Value = 450
$ws.Range("K292").Value = 1000
$ws.Range("L292").Value = 1000
$ws.Range("M292").Value = 1000
$ws.Range("N292").Value = "`$/paquete 20 unidades (volumen en unidades)"
$ws.Range("O292").Value = "Provincia de Talagante"
$ws.Range("P292").Value = 50
$ws.Range("Q292").Value = 20
$ws.Range("R292").Value = "Hortaliza"
